$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5056.7144
$ws.Range("I18").Value = 5498.5
$ws.Range("K18").Value = 5498.5
$ws.Range("M18").Value = -5214.5
$ws.Range("H76").Value = 3474.5
$ws.Range("I76").Value = 2898.75
$ws.Range("J76").Value = 3762.375
$ws.Range("K76").Value = 2898.75
$ws.Range("L76").Value = 3762.375
$ws.Range("M76").Value = -2583.75
$ws.Range("N76").Value = -4392.375
$ws.Range("H79").Value = 3474.5
$ws.Range("I79").Value = 2898.75
$ws.Range("J79").Value = 3762.375
$ws.Range("K79").Value = 2898.75
$ws.Range("L79").Value = 3762.375
$ws.Range("M79").Value = -1806.75
$ws.Range("N79").Value = -5946.375
$ws.Range("H116").Value = 2824.8125
$ws.Range("I116").Value = 2629.7
$ws.Range("K116").Value = 2629.7
$ws.Range("M116").Value = 812.3000000000002
$ws.Range("H135").Value = 4113.15
$ws.Range("I135").Value = 2087.7856
$ws.Range("K135").Value = 18790.0704
$ws.Range("M135").Value = -16255.0704

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1576.9445
$ws.Range("I2").Value = 1086.625
$ws.Range("K2").Value = 1086.625
$ws.Range("M2").Value = -973.625
$ws.Range("H32").Value = 2609.68
$ws.Range("I32").Value = 2598.9321
$ws.Range("J32").Value = 2649.3125
$ws.Range("K32").Value = 2598.9321
$ws.Range("L32").Value = 2649.3125
$ws.Range("M32").Value = -2311.9321
$ws.Range("N32").Value = -3223.3125
$ws.Range("H116").Value = 1576.9445
$ws.Range("I116").Value = 1086.625
$ws.Range("K116").Value = 1086.625
$ws.Range("M116").Value = 1207.375
$ws.Range("H122").Value = 6004.6206
$ws.Range("I122").Value = 5746.6313
$ws.Range("J122").Value = 6494.8
$ws.Range("K122").Value = 17239.8939
$ws.Range("L122").Value = 19484.4
$ws.Range("M122").Value = -14789.8939
$ws.Range("N122").Value = -24384.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1576.9445
$ws.Range("I3").Value = 1086.625
$ws.Range("K3").Value = 1086.625
$ws.Range("M3").Value = -972.625
$ws.Range("H5").Value = 695.8333
$ws.Range("I5").Value = 335
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 335
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = -222
$ws.Range("N5").Value = -2726
$ws.Range("H105").Value = 163.33333
$ws.Range("I105").Value = 163.33333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 163.33333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1583.66667
$ws.Range("N105").ClearContents()
$ws.Range("H131").Value = 24750
$ws.Range("J131").Value = 24750
$ws.Range("L131").Value = 24750
$ws.Range("N131").Value = -34830
$ws.Range("H134").Value = 1733.439
$ws.Range("I134").Value = 1489.7222
$ws.Range("K134").Value = 4469.1666
$ws.Range("M134").Value = -1934.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1854.8125
$ws.Range("I16").Value = 1776.9286
$ws.Range("K16").Value = 1776.9286
$ws.Range("M16").Value = -1489.9286
$ws.Range("H31").Value = 1705.3
$ws.Range("I31").Value = 1233.8
$ws.Range("K31").Value = 1233.8
$ws.Range("M31").Value = -938.8
$ws.Range("H34").Value = 1705.3
$ws.Range("I34").Value = 1233.8
$ws.Range("K34").Value = 1233.8
$ws.Range("M34").Value = -1031.8
$ws.Range("H58").Value = 2594.923
$ws.Range("I58").Value = 2235.5
$ws.Range("K58").Value = 2235.5
$ws.Range("M58").Value = -2032.5
$ws.Range("H99").Value = 3783.25
$ws.Range("I99").Value = 3615.5293
$ws.Range("K99").Value = 3615.5293
$ws.Range("M99").Value = -2117.5293
$ws.Range("H105").Value = 9501.75
$ws.Range("I105").Value = 9501.75
$ws.Range("K105").Value = 9501.75
$ws.Range("M105").Value = -7754.75
$ws.Range("H107").Value = 1097.909
$ws.Range("I107").Value = 725
$ws.Range("K107").Value = 725
$ws.Range("M107").Value = 1195
$ws.Range("H113").Value = 1854.8125
$ws.Range("I113").Value = 1776.9286
$ws.Range("K113").Value = 1776.9286
$ws.Range("M113").Value = 393.0714
$ws.Range("H126").Value = 3783.25
$ws.Range("I126").Value = 3615.5293
$ws.Range("K126").Value = 10846.5879
$ws.Range("M126").Value = -8376.5879
$ws.Range("H130").Value = 133177.8
$ws.Range("J130").Value = 133177.8
$ws.Range("L130").Value = 133177.8
$ws.Range("N130").Value = -143217.8
$ws.Range("H132").Value = 1840.7727
$ws.Range("I132").Value = 1596.4667
$ws.Range("J132").Value = 2364.2856
$ws.Range("K132").Value = 4789.4001
$ws.Range("L132").Value = 7092.8568
$ws.Range("M132").Value = -2259.4001
$ws.Range("N132").Value = -12152.8568
$ws.Range("H134").Value = 1731.375
$ws.Range("I134").Value = 1138.6154
$ws.Range("K134").Value = 3415.8462
$ws.Range("M134").Value = -880.8462
$ws.Range("H136").Value = 2594.923
$ws.Range("I136").Value = 2235.5
$ws.Range("K136").Value = 6706.5
$ws.Range("M136").Value = -4156.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 763.7
$ws.Range("I35").Value = 311.5
$ws.Range("J35").Value = 1065.1666
$ws.Range("K35").Value = 934.5
$ws.Range("L35").Value = 3195.4998
$ws.Range("M35").Value = -646.5
$ws.Range("N35").Value = -3771.4998
$ws.Range("H68").Value = 3102.2222
$ws.Range("I68").Value = 1250
$ws.Range("K68").Value = 3750
$ws.Range("M68").Value = -2939
$ws.Range("H71").Value = 3102.2222
$ws.Range("I71").Value = 1250
$ws.Range("K71").Value = 11250
$ws.Range("M71").Value = -7194
$ws.Range("H93").Value = 4997.5
$ws.Range("I93").Value = 4997.5
$ws.Range("K93").Value = 14992.5
$ws.Range("M93").Value = -13120.5
$ws.Range("H107").Value = 3146.524
$ws.Range("J107").Value = 4289.8
$ws.Range("L107").Value = 12869.4
$ws.Range("N107").Value = -16709.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 41832
$ws.Range("J68").Value = 41832
$ws.Range("L68").Value = 41832
$ws.Range("N68").Value = -43454
$ws.Range("H71").Value = 41832
$ws.Range("J71").Value = 41832
$ws.Range("L71").Value = 125496
$ws.Range("N71").Value = -133608
$ws.Range("H80").Value = 18685
$ws.Range("I80").Value = 18999
$ws.Range("J80").Value = 17900
$ws.Range("K80").Value = 18999
$ws.Range("L80").Value = 17900
$ws.Range("M80").Value = -18001
$ws.Range("N80").Value = -19896
$ws.Range("H83").Value = 18685
$ws.Range("I83").Value = 18999
$ws.Range("J83").Value = 17900
$ws.Range("K83").Value = 94995
$ws.Range("L83").Value = 89500
$ws.Range("M83").Value = -90003
$ws.Range("N83").Value = -99484
$ws.Range("H97").Value = 2347.9167
$ws.Range("I97").Value = 1275.5625
$ws.Range("J97").Value = 4492.625
$ws.Range("K97").Value = 1275.5625
$ws.Range("L97").Value = 4492.625
$ws.Range("M97").Value = -779.5625
$ws.Range("N97").Value = -5484.625
$ws.Range("H102").Value = 50994
$ws.Range("I102").Value = 935.6667
$ws.Range("K102").Value = 935.6667
$ws.Range("M102").Value = 686.3333
$ws.Range("H109").Value = 84443
$ws.Range("J109").Value = 84443
$ws.Range("L109").Value = 84443
$ws.Range("N109").Value = -86523
$ws.Range("H132").Value = 3790.68
$ws.Range("I132").Value = 3314.6826
$ws.Range("K132").Value = 9944.0478
$ws.Range("M132").Value = -7414.0478

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9733.514999999999
$ws.Range("I61").Value = 8716.24
$ws.Range("J61").Value = 12276.7
$ws.Range("K61").Value = 8716.24
$ws.Range("L61").Value = 12276.7
$ws.Range("M61").Value = -8514.24
$ws.Range("N61").Value = -12680.7
$ws.Range("H113").Value = 9733.514999999999
$ws.Range("I113").Value = 8716.24
$ws.Range("J113").Value = 12276.7
$ws.Range("K113").Value = 8716.24
$ws.Range("L113").Value = 12276.7
$ws.Range("M113").Value = -6546.24
$ws.Range("N113").Value = -16616.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 263059
$ws.Range("I62").Value = 263059
$ws.Range("K62").Value = 263059
$ws.Range("M62").Value = -262435
$ws.Range("H65").Value = 263059
$ws.Range("I65").Value = 263059
$ws.Range("K65").Value = 1315295
$ws.Range("M65").Value = -1312175
$ws.Range("H100").Value = 1875.4286
$ws.Range("I100").Value = 892.5625
$ws.Range("J100").Value = 5020.6
$ws.Range("K100").Value = 1785.125
$ws.Range("L100").Value = 10041.2
$ws.Range("M100").Value = -1244.125
$ws.Range("N100").Value = -11123.2
$ws.Range("H107").Value = 741.0345
$ws.Range("I107").Value = 476.35294
$ws.Range("K107").Value = 1429.05882
$ws.Range("M107").Value = 490.94118
$ws.Range("H113").Value = 3789016
$ws.Range("I113").Value = 6945264
$ws.Range("K113").Value = 20835792
$ws.Range("M113").Value = -20833622
$ws.Range("H122").Value = 7835.4443
$ws.Range("J122").Value = 6902.1665
$ws.Range("L122").Value = 20706.4995
$ws.Range("N122").Value = -25606.4995
